# Added Jan 28 Tournament: appends 14 new match rows (62-75) to the Elo time table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell whose style (bold/centered/bordered) is reused for the new
# "match index" cells in column A.
$styleSource = $ws.Range("A61")

# Row 62: match 61
$styleSource.Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "'"
$ws.Range("B62").Style = "Normal"
$ws.Range("C62").Value = -17.40927553185213
$ws.Range("D62").Value = "'"
$ws.Range("D62").Style = "Normal"
$ws.Range("E62").Value = "'"
$ws.Range("E62").Style = "Normal"
$ws.Range("F62").Value = 17.40927553185213

# Row 63: match 62
$styleSource.Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 6.146545249736358
$ws.Range("C63").Value = "'"
$ws.Range("C63").Style = "Normal"
$ws.Range("D63").Value = "'"
$ws.Range("D63").Style = "Normal"
$ws.Range("E63").Value = -6.146545249736356
$ws.Range("F63").Value = "'"
$ws.Range("F63").Style = "Normal"

# Row 64: match 63
$styleSource.Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = -23.03791779652222
$ws.Range("C64").Value = 23.03791779652222
$ws.Range("D64").Value = "'"
$ws.Range("D64").Style = "Normal"
$ws.Range("E64").Value = "'"
$ws.Range("E64").Style = "Normal"
$ws.Range("F64").Value = "'"
$ws.Range("F64").Style = "Normal"

# Row 65: match 64
$styleSource.Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = "'"
$ws.Range("B65").Style = "Normal"
$ws.Range("C65").Value = "'"
$ws.Range("C65").Style = "Normal"
$ws.Range("D65").Value = "'"
$ws.Range("D65").Style = "Normal"
$ws.Range("E65").Value = -10.95195006581875
$ws.Range("F65").Value = 10.95195006581875

# Row 66: match 65
$styleSource.Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = "'"
$ws.Range("B66").Style = "Normal"
$ws.Range("C66").Value = 9.327181778093486
$ws.Range("D66").Value = "'"
$ws.Range("D66").Style = "Normal"
$ws.Range("E66").Value = -9.327181778093484
$ws.Range("F66").Value = "'"
$ws.Range("F66").Style = "Normal"

# Row 67: match 66
$styleSource.Copy()
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 11.3121180763797
$ws.Range("C67").Value = "'"
$ws.Range("C67").Style = "Normal"
$ws.Range("D67").Value = "'"
$ws.Range("D67").Style = "Normal"
$ws.Range("E67").Value = "'"
$ws.Range("E67").Style = "Normal"
$ws.Range("F67").Value = -11.3121180763797

# Row 68: match 67
$styleSource.Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "'"
$ws.Range("B68").Style = "Normal"
$ws.Range("C68").Value = -17.23351129782776
$ws.Range("D68").Value = "'"
$ws.Range("D68").Style = "Normal"
$ws.Range("E68").Value = "'"
$ws.Range("E68").Style = "Normal"
$ws.Range("F68").Value = 17.23351129782776

# Row 69: match 68
$styleSource.Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = 5.095205994447398
$ws.Range("C69").Value = "'"
$ws.Range("C69").Style = "Normal"
$ws.Range("D69").Value = "'"
$ws.Range("D69").Style = "Normal"
$ws.Range("E69").Value = -5.095205994447396
$ws.Range("F69").Value = "'"
$ws.Range("F69").Style = "Normal"

# Row 70: match 69
$styleSource.Copy()
$ws.Range("A70").PasteSpecial(-4122)
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = 9.456519108265985
$ws.Range("C70").Value = -9.456519108265981
$ws.Range("D70").Value = "'"
$ws.Range("D70").Style = "Normal"
$ws.Range("E70").Value = "'"
$ws.Range("E70").Style = "Normal"
$ws.Range("F70").Value = "'"
$ws.Range("F70").Style = "Normal"

# Row 71: match 70
$styleSource.Copy()
$ws.Range("A71").PasteSpecial(-4122)
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "'"
$ws.Range("B71").Style = "Normal"
$ws.Range("C71").Value = "'"
$ws.Range("C71").Style = "Normal"
$ws.Range("D71").Value = "'"
$ws.Range("D71").Style = "Normal"
$ws.Range("E71").Value = -7.924851840315794
$ws.Range("F71").Value = 7.924851840315796

# Row 72: match 71
$styleSource.Copy()
$ws.Range("A72").PasteSpecial(-4122)
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "'"
$ws.Range("B72").Style = "Normal"
$ws.Range("C72").Value = 8.974044534140853
$ws.Range("D72").Value = "'"
$ws.Range("D72").Style = "Normal"
$ws.Range("E72").Value = -8.974044534140855
$ws.Range("F72").Value = "'"
$ws.Range("F72").Style = "Normal"

# Row 73: match 72
$styleSource.Copy()
$ws.Range("A73").PasteSpecial(-4122)
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = -20.61844630140414
$ws.Range("C73").Value = "'"
$ws.Range("C73").Style = "Normal"
$ws.Range("D73").Value = "'"
$ws.Range("D73").Style = "Normal"
$ws.Range("E73").Value = "'"
$ws.Range("E73").Style = "Normal"
$ws.Range("F73").Value = 20.61844630140414

# Row 74: match 73
$styleSource.Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = -17.2985575863854
$ws.Range("D74").Value = "'"
$ws.Range("D74").Style = "Normal"
$ws.Range("E74").Value = "'"
$ws.Range("E74").Style = "Normal"
$ws.Range("F74").Value = 17.2985575863854

# Row 75: match 74
$styleSource.Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = -14.36305800477309
$ws.Range("D75").Value = "'"
$ws.Range("D75").Style = "Normal"
$ws.Range("E75").Value = "'"
$ws.Range("E75").Style = "Normal"
$ws.Range("F75").Value = 14.36305800477309
